# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates to the Leviathan_Profits workbook sheets
# as captured by the authoritative OOXML diff (scheduled market-data refresh).

$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3359.1667
$ws.Range("J32").Value = 3231
$ws.Range("L32").Value = 3231
$ws.Range("N32").Value = -3883
$ws.Range("H100").Value = 8159.0835
$ws.Range("I100").Value = 8390.9
$ws.Range("K100").Value = 8390.9
$ws.Range("M100").Value = -7849.9
$ws.Range("H113").Value = 41739.81
$ws.Range("I113").Value = 93157.17999999999
$ws.Range("J113").Value = 4033.7334
$ws.Range("K113").Value = 93157.17999999999
$ws.Range("L113").Value = 4033.7334
$ws.Range("M113").Value = -89903.17999999999
$ws.Range("N113").Value = -10541.7334
$ws.Range("H116").Value = 7522.4194
$ws.Range("I116").Value = 9057.1
$ws.Range("K116").Value = 9057.1
$ws.Range("M116").Value = -5615.1

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1962.5
$ws.Range("I61").Value = 2045
$ws.Range("K61").Value = 2045
$ws.Range("M61").Value = -1833
$ws.Range("H74").Value = 2021.5
$ws.Range("J74").Value = 1890.6666
$ws.Range("L74").Value = 1890.6666
$ws.Range("N74").Value = -3638.6666
$ws.Range("H77").Value = 2021.5
$ws.Range("J77").Value = 1890.6666
$ws.Range("L77").Value = 9453.333000000001
$ws.Range("N77").Value = -18189.333
$ws.Range("H97").Value = 1026.3636
$ws.Range("I97").Value = 1085.5862
$ws.Range("K97").Value = 1085.5862
$ws.Range("M97").Value = -589.5862
$ws.Range("H102").Value = 1421.3448
$ws.Range("I102").Value = 1387.32
$ws.Range("J102").Value = 1634
$ws.Range("K102").Value = 1387.32
$ws.Range("L102").Value = 1634
$ws.Range("M102").Value = 234.6800000000001
$ws.Range("N102").Value = -4878
$ws.Range("H132").Value = 1616.0952
$ws.Range("I132").Value = 1394.7222
$ws.Range("J132").Value = 2944.3333
$ws.Range("K132").Value = 4184.1666
$ws.Range("L132").Value = 8832.999899999999
$ws.Range("M132").Value = -1654.1666
$ws.Range("N132").Value = -13892.9999
$ws.Range("H133").Value = 99126.5
$ws.Range("J133").Value = 99126.5
$ws.Range("L133").Value = 99126.5
$ws.Range("N133").Value = -104186.5
$ws.Range("H136").Value = 1962.5
$ws.Range("I136").Value = 2045
$ws.Range("K136").Value = 6135
$ws.Range("M136").Value = -3585

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1845
$ws.Range("J94").Value = 1499
$ws.Range("L94").Value = 1499
$ws.Range("N94").Value = -2401
$ws.Range("H99").Value = 2001.25
$ws.Range("I99").Value = 2001.25
$ws.Range("K99").Value = 2001.25
$ws.Range("M99").Value = -503.25
$ws.Range("H105").Value = 11141.308
$ws.Range("I105").Value = 13173.7
$ws.Range("K105").Value = 13173.7
$ws.Range("M105").Value = -11426.7

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32220.045
$ws.Range("I31").Value = 30310.027
$ws.Range("K31").Value = 30310.027
$ws.Range("M31").Value = -30015.027
$ws.Range("H34").Value = 32220.045
$ws.Range("I34").Value = 30310.027
$ws.Range("K34").Value = 30310.027
$ws.Range("M34").Value = -30108.027
$ws.Range("H62").Value = 2912.4285
$ws.Range("I62").Value = 2958.6
$ws.Range("K62").Value = 2958.6
$ws.Range("M62").Value = -2334.6
$ws.Range("H65").Value = 2912.4285
$ws.Range("I65").Value = 2958.6
$ws.Range("K65").Value = 14793
$ws.Range("M65").Value = -11673
$ws.Range("H97").Value = 36899
$ws.Range("J97").Value = 36899
$ws.Range("L97").Value = 36899
$ws.Range("N97").Value = -38881
$ws.Range("H105").Value = 1346.8
$ws.Range("I105").Value = 907.6667
$ws.Range("K105").Value = 907.6667
$ws.Range("M105").Value = 839.3333
$ws.Range("H122").Value = 84409.664
$ws.Range("I122").Value = 111991.22
$ws.Range("J122").Value = 1665
$ws.Range("K122").Value = 335973.66
$ws.Range("L122").Value = 4995
$ws.Range("M122").Value = -333523.66
$ws.Range("N122").Value = -9895

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 162.25
$ws.Range("I12").Value = 107.75
$ws.Range("J12").Value = 189.5
$ws.Range("K12").Value = 323.25
$ws.Range("L12").Value = 568.5
$ws.Range("M12").Value = -150.25
$ws.Range("N12").Value = -914.5
$ws.Range("H23").Value = 112.181816
$ws.Range("J23").Value = 86
$ws.Range("L23").Value = 258
$ws.Range("N23").Value = -728
$ws.Range("H64").Value = 1649.75
$ws.Range("J64").Value = 1933
$ws.Range("L64").Value = 5799
$ws.Range("N64").Value = -6339
$ws.Range("H67").Value = 1649.75
$ws.Range("J67").Value = 1933
$ws.Range("L67").Value = 5799
$ws.Range("N67").Value = -7671
$ws.Range("H80").Value = 2400
$ws.Range("J80").Value = 2400
$ws.Range("L80").Value = 7200
$ws.Range("N80").Value = -9072
$ws.Range("H83").Value = 2400
$ws.Range("J83").Value = 2400
$ws.Range("L83").Value = 21600
$ws.Range("N83").Value = -30960
$ws.Range("H98").Value = 1071
$ws.Range("I98").Value = 761.3333
$ws.Range("K98").Value = 2283.9999
$ws.Range("M98").Value = -785.9998999999998
$ws.Range("H104").Value = 4423.625
$ws.Range("I104").Value = 2699.5
$ws.Range("J104").Value = 4998.3335
$ws.Range("K104").Value = 8098.5
$ws.Range("L104").Value = 14995.0005
$ws.Range("M104").Value = -5477.5
$ws.Range("N104").Value = -20237.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 55555
$ws.Range("J24").Value = 55555
$ws.Range("L24").Value = 55555
$ws.Range("N24").Value = -55901
$ws.Range("H70").Value = 8433.272000000001
$ws.Range("I70").Value = 8307.888999999999
$ws.Range("J70").Value = 8997.5
$ws.Range("K70").Value = 8307.888999999999
$ws.Range("L70").Value = 8997.5
$ws.Range("M70").Value = -8037.888999999999
$ws.Range("N70").Value = -9537.5
$ws.Range("H73").Value = 8433.272000000001
$ws.Range("I73").Value = 8307.888999999999
$ws.Range("J73").Value = 8997.5
$ws.Range("K73").Value = 8307.888999999999
$ws.Range("L73").Value = 8997.5
$ws.Range("M73").Value = -7371.888999999999
$ws.Range("N73").Value = -10869.5
$ws.Range("H126").Value = 2678.5454
$ws.Range("I126").Value = 2575.5715
$ws.Range("J126").Value = 2858.75
$ws.Range("K126").Value = 7726.7145
$ws.Range("L126").Value = 8576.25
$ws.Range("M126").Value = -5256.7145
$ws.Range("N126").Value = -13516.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 38421.89
$ws.Range("I7").Value = 53633
$ws.Range("K7").Value = 53633
$ws.Range("M7").Value = -53521
$ws.Range("H16").Value = 180639.9
$ws.Range("I16").Value = 115056.71
$ws.Range("K16").Value = 115056.71
$ws.Range("M16").Value = -114886.71
$ws.Range("H69").Value = 46250
$ws.Range("J69").Value = 46250
$ws.Range("L69").Value = 46250
$ws.Range("N69").Value = -47872
$ws.Range("H72").Value = 46250
$ws.Range("J72").Value = 46250
$ws.Range("L72").Value = 138750
$ws.Range("N72").Value = -146862
$ws.Range("H93").Value = 9453.15
$ws.Range("I93").Value = 1171.5526
$ws.Range("K93").Value = 1171.5526
$ws.Range("M93").Value = 76.44740000000002
$ws.Range("H97").Value = 17500
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H126").Value = 38421.89
$ws.Range("I126").Value = 53633
$ws.Range("K126").Value = 160899
$ws.Range("M126").Value = -158429
$ws.Range("H133").Value = 107498.75
$ws.Range("J133").Value = 107498.75
$ws.Range("L133").Value = 107498.75
$ws.Range("N133").Value = -112558.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 5477.222
$ws.Range("J51").Value = 40495
$ws.Range("L51").Value = 40495
$ws.Range("N51").Value = -41515
$ws.Range("H52").Value = 27499.5
$ws.Range("I52").Value = 20000
$ws.Range("J52").Value = 34999
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 34999
$ws.Range("M52").Value = -19774
$ws.Range("N52").Value = -35451
$ws.Range("H100").Value = 2240.5518
$ws.Range("I100").Value = 2132.8572
$ws.Range("K100").Value = 4265.7144
$ws.Range("M100").Value = -3724.7144
$ws.Range("H122").Value = 1706.6666
$ws.Range("I122").Value = 1622.8572
$ws.Range("K122").Value = 4868.571599999999
$ws.Range("M122").Value = -2418.571599999999
